$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: min_samples_leaf
$ws.Range("A2").Value = "min_samples_leaf"
$ws.Range("B2").Value = 1045.296
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 9.464
$ws.Range("E2").Value = 0.005955801935382763

# Row 3: min_weight_fraction_leaf
$ws.Range("A3").Value = "min_weight_fraction_leaf"
$ws.Range("B3").Value = 33820.942
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 306.206
$ws.Range("E3").Value = [double]"1.360934806132981E-13"

# Row 4: max_features
$ws.Range("A4").Value = "max_features"
$ws.Range("B4").Value = 3509.589
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 31.775
$ws.Range("E4").Value = [double]"1.619632069649968E-05"

# Row 5: min_samples_leaf:min_weight_fraction_leaf
$ws.Range("A5").Value = "min_samples_leaf:min_weight_fraction_leaf"
$ws.Range("B5").Value = 1494.239
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 13.528
$ws.Range("E5").Value = 0.001491337832162908

# Row 6: min_samples_leaf:max_features
$ws.Range("A6").Value = "min_samples_leaf:max_features"
$ws.Range("B6").Value = 0.917
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0.008304142741944997
$ws.Range("E6").Value = 0.928

# Row 7: min_weight_fraction_leaf:max_features
$ws.Range("A7").Value = "min_weight_fraction_leaf:max_features"
$ws.Range("B7").Value = 1858.979
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 16.831
$ws.Range("E7").Value = 0.0005535449846343956

# Row 8: Residual
$ws.Range("A8").Value = "Residual"
$ws.Range("B8").Value = 2209.032
$ws.Range("C8").Value = 20
$ws.Range("D8").Value = $null
$ws.Range("E8").Value = $null

# Remove old rows 9-12 entirely (data no longer present)
$ws.Range("A9:E12").Clear()
